# Update "want to go" counts (column F) on several sheets to reflect the
# latest generated output (gh-pages build at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 109
$ws1.Range("F4").Value = 1539
$ws1.Range("F5").Value = 234
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 671
$ws1.Range("F8").Value = 10037
$ws1.Range("F10").Value = 127
$ws1.Range("F11").Value = 246
$ws1.Range("F12").Value = 189
$ws1.Range("F13").Value = 380
$ws1.Range("F14").Value = 6955
$ws1.Range("F16").Value = 650
$ws1.Range("F18").Value = 211

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 6
$ws2.Range("F3").Value = 553

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 109
$ws4.Range("F4").Value = 1539
$ws4.Range("F5").Value = 234
$ws4.Range("F6").Value = 6
$ws4.Range("F7").Value = 50
$ws4.Range("F8").Value = 671
$ws4.Range("F9").Value = 553
$ws4.Range("F11").Value = 10037
$ws4.Range("F13").Value = 127
$ws4.Range("F14").Value = 246
$ws4.Range("F15").Value = 189
$ws4.Range("F16").Value = 380
$ws4.Range("F17").Value = 6955
$ws4.Range("F19").Value = 650
$ws4.Range("F21").Value = 211

$wb.Save()
